# Auto-generated script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.844.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.851.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.623.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.837.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.389.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.16%  "
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.762.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
